# Update mortality tracker data: add rows for 18-21 April 2024.
#
# The WayBackMachine failed to save the Aljazeera tracking site correctly on
# 21 April 2024; the casualty figures recorded for 18-20 April are repeats of
# data already captured, and the 21 April row carries the newest numbers but
# (per the commit message) has no working archive.org source link.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 194: tracker_date 18.04.2024, report_date carried over from 17.04.2024 ---
$ws.Range("A194").Value = "18.04.2024"
$ws.Range("B194").Value = "17.04.2024"
$ws.Range("C194").Value = 33899
$ws.Range("D194").Value = 13800
$ws.Range("E194").Value = 8400
$ws.Range("F194").Value = 76664
$ws.Range("I194").Value = 8000
$ws.Range("J194").Value = 468
$ws.Range("K194").Value = 120
$ws.Range("L194").Value = 4750
$ws.Range("M194").Value = "https://web.archive.org/web/20240418215003/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"
$ws.Range("C194").WrapText = $true
$ws.Range("F194").WrapText = $true

# --- Row 195: tracker_date 19.04.2024, report_date 19.04.2024 ---
$ws.Range("A195").Value = "19.04.2024"
$ws.Range("B195").Value = "19.04.2024"
$ws.Range("C195").Value = 34012
$ws.Range("D195").Value = 13800
$ws.Range("E195").Value = 8400
$ws.Range("F195").Value = 76833
$ws.Range("I195").Value = 8000
$ws.Range("J195").Value = 468
$ws.Range("K195").Value = 120
$ws.Range("L195").Value = 4800
$ws.Range("M195").Value = "https://web.archive.org/web/20240419203633/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"

# --- Row 196: tracker_date 20.04.2024, report_date still 19.04.2024 (archive glitch) ---
$ws.Range("A196").Value = "20.04.2024"
$ws.Range("B196").Value = "19.04.2024"
$ws.Range("C196").Value = 34012
$ws.Range("D196").Value = 13800
$ws.Range("E196").Value = 8400
$ws.Range("F196").Value = 76833
$ws.Range("I196").Value = 8000
$ws.Range("J196").Value = 468
$ws.Range("K196").Value = 120
$ws.Range("L196").Value = 4800
$ws.Range("M196").Value = "https://web.archive.org/web/20240420210539/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"

# --- Row 197: tracker_date 21.04.2024, report_date 21.04.2024, no working source link ---
$ws.Range("A197").Value = "21.04.2024"
$ws.Range("B197").Value = "21.04.2024"
$ws.Range("C197").Value = 34097
$ws.Range("D197").Value = 13800
$ws.Range("E197").Value = 8400
$ws.Range("F197").Value = 76980
$ws.Range("I197").Value = 8000
$ws.Range("J197").Value = 485
$ws.Range("K197").Value = 122
$ws.Range("L197").Value = 4800
$ws.Range("C197").WrapText = $true
$ws.Range("F197").WrapText = $true

# Leave the cursor where the editor's session ended up, at the newest entry.
$ws.Range("M204").Select()
